$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the totals row (row 31), shifting the totals
# row and the footer row down by one. This mirrors Excel's row-insert
# behaviour (existing merged cells below the insertion point shift too).
$ws.Rows("31:31").Insert()

# Clone the formatting of the row above (the last product row) onto the
# newly inserted row so the new row gets the same cell styles used by
# every other product row in the table.
$ws.Range("A30:Q30").Copy()
$ws.Range("A31:Q31").PasteSpecial(-4122)

# New merged regions for row 31 matching the other product rows.
$ws.Range("A31:B31").Merge()
$ws.Range("C31:G31").Merge()
$ws.Range("H31:K31").Merge()
$ws.Range("L31:M31").Merge()
$ws.Range("N31:O31").Merge()

# Populate the new product row (#25). A31/C31/H31/N31/Q31 are already
# formatted as Text (numFmtId 49) after the format copy, so plain string
# assignment stores them as text. L31/P31 carry a numeric display format
# (qty-ratio / 2-decimal), so a bare string assignment would be coerced
# to a real number there; temporarily flipping those two to Text for the
# assignment (then restoring the original number format) keeps the cell
# genuinely text -- matching how this report stores those two columns --
# without leaving the cell's number format changed.
$ws.Range("A31").Value = 25
$ws.Range("C31").Value = "كريم فيرند لافلي الصغير"
$ws.Range("H31").Value = "4:0"

$fmtL31 = $ws.Range("L31").NumberFormat()
$ws.Range("L31").NumberFormat = "@"
$ws.Range("L31").Value = "0"
$ws.Range("L31").NumberFormat = $fmtL31

$ws.Range("N31").Value = "20.00"

$fmtP31 = $ws.Range("P31").NumberFormat()
$ws.Range("P31").NumberFormat = "@"
$ws.Range("P31").Value = "20.0000"
$ws.Range("P31").NumberFormat = $fmtP31

$ws.Range("Q31").Value = "1/1"

# Update the grand-total cell (now shifted to row 32).
$ws.Range("P32").Value = 1356.75

# Refresh the generated-on timestamp in the footer (now on row 33).
$ws.Range("A33").Value = "Tuesday, 2 September, 2025 11:28 AM"
